$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" column (E16:E21) listed the account-statement years in
# descending order (2009..2004). Update the base/database so the years are
# listed in ascending order (2004..2009) instead.
$ws.Range("E16").Value = 2004
$ws.Range("E17").Value = 2005
$ws.Range("E18").Value = 2006
$ws.Range("E19").Value = 2007
$ws.Range("E20").Value = 2008
$ws.Range("E21").Value = 2009
